# Auto-generated edit script
# Applies text replacements to answers-of-addition_and_subtraction_within_100.docx
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2026-01-06 Tuesday" "2026-01-07 Wednesday"
Replace-Text "77+6=83" "69-54=15"
Replace-Text "46+47=93" "40+21=61"
Replace-Text "79-0=79" "25+60=85"
Replace-Text "74-36=38" "61-38=23"
Replace-Text "85-84=1" "29+50=79"
Replace-Text "5+52=57" "77-31=46"
Replace-Text "33+8=41" "25+37=62"
Replace-Text "50-10=40" "3+63=66"
Replace-Text "20+65=85" "71-1=70"
Replace-Text "76-29=47" "52-12=40"
Replace-Text "96-45=51" "89-70=19"
Replace-Text "55-3=52" "54-15=39"
Replace-Text "47-37=10" "33+33=66"
Replace-Text "23+5=28" "21+11=32"
Replace-Text "53+5=58" "41+2=43"
Replace-Text "47+2=49" "76-50=26"
Replace-Text "48-10=38" "18+29=47"
Replace-Text "79-51=28" "9+55=64"
Replace-Text "53-42=11" "32+25=57"
Replace-Text "21+15=36" "59+33=92"
Replace-Text "7+60=67" "49+30=79"
Replace-Text "46+11=57" "42+26=68"
Replace-Text "93-77=16" "51+0=51"
Replace-Text "31+52=83" "40+11=51"
Replace-Text "84-1=83" "53-0=53"
Replace-Text "23+60=83" "96-7=89"
Replace-Text "40-4=36" "42-30=12"
Replace-Text "90-63=27" "86-1=85"
Replace-Text "35+32=67" "72+1=73"
Replace-Text "50-20=30" "73-32=41"
Replace-Text "49+7=56" "88-32=56"
Replace-Text "87-66=21" "68+12=80"
Replace-Text "88-51=37" "82+3=85"
Replace-Text "64+20=84" "6+51=57"
Replace-Text "13+68=81" "70+14=84"
Replace-Text "5+46=51" "93-87=6"
Replace-Text "15+38=53" "53-8=45"
Replace-Text "15+83=98" "28-7=21"
Replace-Text "72+21=93" "68-31=37"
Replace-Text "66-18=48" "2-0=2"
Replace-Text "29-0=29" "14+25=39"
Replace-Text "91-64=27" "89+2=91"
Replace-Text "42+8=50" "85-42=43"
Replace-Text "22+56=78" "53-25=28"
Replace-Text "41-9=32" "66+27=93"
Replace-Text "21-18=3" "31+19=50"
Replace-Text "52+35=87" "33+25=58"
Replace-Text "46-7=39" "30+39=69"
Replace-Text "85-21=64" "47-6=41"
Replace-Text "54-5=49" "47+7=54"
Replace-Text "93+3=96" "8+57=65"
Replace-Text "53+8=61" "20-12=8"
Replace-Text "66-30=36" "29+19=48"
Replace-Text "96-59=37" "72-18=54"
Replace-Text "81-41=40" "38+7=45"
Replace-Text "99-37=62" "77-6=71"
Replace-Text "62-40=22" "18+4=22"
Replace-Text "98-6=92" "56-51=5"
Replace-Text "53-1=52" "33+15=48"
Replace-Text "95-63=32" "30-11=19"
Replace-Text "11+65=76" "17+50=67"
Replace-Text "19+21=40" "16-15=1"
Replace-Text "46+43=89" "19+8=27"
Replace-Text "66-19=47" "19-19=0"
Replace-Text "8+59=67" "77-2=75"
Replace-Text "97-56=41" "22-20=2"
Replace-Text "72-48=24" "22+57=79"
Replace-Text "9+68=77" "91-35=56"
Replace-Text "24-1=23" "41+18=59"
Replace-Text "77-22=55" "35-27=8"
Replace-Text "4+36=40" "17+20=37"
Replace-Text "24+45=69" "29+70=99"
Replace-Text "57-12=45" "40-35=5"
Replace-Text "80-52=28" "59-30=29"
Replace-Text "2+58=60" "0+79=79"
Replace-Text "95-76=19" "82-31=51"
Replace-Text "72-3=69" "43+52=95"
Replace-Text "4+28=32" "6+6=12"
Replace-Text "83-26=57" "85-36=49"
Replace-Text "16+11=27" "9+22=31"
Replace-Text "13+22=35" "66-18=48"
Replace-Text "59+13=72" "67-24=43"
Replace-Text "85-70=15" "25-12=13"
Replace-Text "40+35=75" "81-31=50"
Replace-Text "30+41=71" "52-5=47"
Replace-Text "79-77=2" "7+68=75"
Replace-Text "24+17=41" "99-40=59"
Replace-Text "2+66=68" "20+72=92"
Replace-Text "91-69=22" "23+46=69"
Replace-Text "15+19=34" "45-13=32"
Replace-Text "20+75=95" "15+77=92"
Replace-Text "41+31=72" "78-3=75"
Replace-Text "30-25=5" "14+53=67"
Replace-Text "74-47=27" "44+9=53"
Replace-Text "56+22=78" "55+24=79"
Replace-Text "2+60=62" "13-5=8"
Replace-Text "28+55=83" "28+26=54"
Replace-Text "93-73=20" "41+50=91"
Replace-Text "96-69=27" "83-6=77"
Replace-Text "28-18=10" "71-46=25"
